$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Remove the now-unused Sheet2 and Sheet3 ---
[void]$wb.Worksheets.Item("Sheet2").Delete()
[void]$wb.Worksheets.Item("Sheet3").Delete()

$ws = $wb.Worksheets.Item("Sheet1")

# --- Update the status text per-row (new shared strings, in desired final order) ---
$ws.Range("B4").Value = "la timp"
$ws.Range("B5").Value = "la timp"
$ws.Range("B6").Value = "la timp"
$ws.Range("B7").Value = "sosește la timp la Târgoviște*"
$ws.Range("B8").Value = "pleacă la timp din Titu*"
$ws.Range("B9").Value = "pleacă la timp din București Nord*"
$ws.Range("B10").Value = "pleacă la timp din București Nord*"
$ws.Range("B11").Value = "sosește la timp la București Nord*"
$ws.Range("B12").Value = "pleacă la timp din Constanța*"
